# Update Name of Algo
# Apply updated RandomForest imputation results to columns A and B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A4"  = -19.78929999999998
    "B4"  = 8.787599999999996
    "A6"  = -22.67920000000002
    "A7"  = -20.07389999999998
    "B9"  = 6.184200000000001
    "B12" = 4.719899999999996
    "A16" = -21.94860000000001
    "B17" = 5.827199999999998
    "B18" = 6.525099999999996
    "B19" = 9.150199999999998
    "A20" = -22.67750000000001
    "B20" = 5.327999999999996
    "B26" = 4.419100000000004
    "A28" = -21.9995
    "A29" = -21.26539999999997
    "B31" = 5.092899999999998
    "A32" = -21.16360000000001
    "B39" = 9.538300000000007
    "A40" = -19.90429999999997
    "B40" = 8.847999999999995
    "B41" = 9.228599999999989
    "B42" = 9.218099999999994
    "B43" = 6.278800000000003
    "A46" = -21.90820000000001
    "B47" = 6.009499999999999
    "B48" = 5.548100000000004
    "A51" = -21.74329999999999
    "A52" = -22.1179
    "A57" = -21.88970000000003
    "A59" = -22.17010000000001
    "A62" = -22.08670000000001
    "B63" = 4.811599999999998
    "B64" = 5.371000000000001
    "A66" = -21.49750000000001
    "A73" = -20.10219999999998
    "A74" = -21.69299999999998
    "B76" = 5.282600000000005
    "B81" = 5.369700000000003
    "B89" = 5.3926
    "A92" = -21.64400000000001
    "B94" = 4.618399999999991
    "A100" = -22.1048
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
